# Apply the "Add files via upload" update to the AAVE/polygon tracker sheet:
#   - D2 (the "combien de dollars j'ajoute a chaque achat" setting) goes from 2 to 3
#   - A4 (first price entry) is corrected from 295.46 to 2
#   - two new purchase rows (19 & 20) are appended with their date/time stamps
#   - the active selection moves to D2
#   - the recalculated summary cells (I2,J2,K2,L2,M2) update automatically via formulas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dollars added per purchase" setting.
$ws.Range("D2").Value = 3

# Correct the first data row's price.
$ws.Range("A4").Value = 2

# Rows 17 & 18 switch from an automatic row height to an explicit custom
# height (same 12.75 value), matching the newly-appended rows below them.
$ws.Rows(17).RowHeight = 12.75
$ws.Rows(18).RowHeight = 12.75

# Append the two new purchase rows.
$ws.Range("A19").Value = 162.34
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "20/11/2025"
$ws.Range("D19").Value = "19:00:16"

$ws.Range("A20").Value = 161.11
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "20/11/2025"
$ws.Range("D20").Value = "20:14:09"

# Move the active cell/selection to D2.
$ws.Range("D2").Select()
